$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rework the "Contact" sheet: remove the old 5-column contact record and
#    replace it with a single "HL Employee" column (header + 4 names).
# ---------------------------------------------------------------------------
$contact = $wb.Worksheets.Item("Contact")

# Drop the mailto: hyperlink before the underlying cell content is cleared.
$contact.Hyperlinks.Delete()

# Clear out the old CompanyName/FirstName/LastName/Email/Phone columns.
$contact.Range("B1:E2").Clear()

# New single-column content.
$contact.Range("A1").Value = "HL Employee"
$contact.Range("A2").Value = "Aaron Solomon"
$contact.Range("A3").Value = "Andy Lund"
$contact.Range("A4").Value = "Jack Truett"
$contact.Range("A5").Value = "Mark Francis"

$contact.Range("A1").Style = "Normal"
$contact.Range("A1").Font.Bold = $true

$contact.Range("D9").Select() | Out-Null

# The Hyperlink cell style is no longer used anywhere in the workbook.
$wb.Styles.Item("Hyperlink").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. Append a brand new "Industry Group" sheet after "Contact".
# ---------------------------------------------------------------------------
$industry = $wb.Worksheets.Add($null, $contact)
$industry.Name = "Industry Group"

$industry.Range("A1").Value = "Industry Group"
$industry.Range("A1").Font.Bold = $true
$industry.Range("A1").HorizontalAlignment = -4108
$industry.Range("A1").VerticalAlignment = -4108

$industry.Range("A2").Value = "--None--"
$industry.Range("A2").NumberFormat = "@"

$industry.Range("A3").Value = "FT - FinTech"

$industry.Columns("A:A").AutoFit() | Out-Null

$industry.Range("B3").Select() | Out-Null

# "Contact" remains the active tab (matches activeTab="1" in bookViews).
$contact.Activate()
